# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns for
# rows 2-51 on Sheet1 to reflect the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.948.06"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.878.80"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7413"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.71"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07216"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08342"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").Value = "2.033.26"
$ws.Range("E12").Value = "  +9.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7528"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.398"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.35"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.143"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "29.983.67"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "248.73"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007857"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "2.134.95"
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.009"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.303"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.10"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.497"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.593"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.537"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.215"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05371"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7534"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01966"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.758"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "1.122.60"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.161"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.66"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8618"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.79"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.868"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.622"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.529"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "2.037.14"
$ws.Range("E51").Value = "  +2.26%  "
